# Generate Report for Handoff
# Refresh the "Latest Handoff Datetime" for the files that were just
# (re-)handed off: 6cbafaf2-..., 0786a1aa-..., e7bb0446-...
# Row 7 = 6cbafaf2-1527-4d17-8609-c520de4665c1
# Row 11 = 0786a1aa-dfac-4cdd-a92e-f3d084e78018
# Row 16 = e7bb0446-1734-43a6-8b4a-14d3e5e8503e

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("D7").Value = "2016-03-09 04:36:23"
$zhcn.Range("D10").Value = "2016-03-09 04:36:23"
$zhcn.Range("D11").Value = "2016-03-09 04:36:23"
$zhcn.Range("D12").Value = "2016-03-09 04:36:23"
$zhcn.Range("D13").Value = "2016-03-09 04:36:23"
$zhcn.Range("D14").Value = "2016-03-09 04:36:23"
$zhcn.Range("D15").Value = "2016-03-09 04:36:23"
$zhcn.Range("D16").Value = "2016-03-09 04:36:23"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("D7").Value = "2016-03-09 04:36:26"
$dede.Range("D10").Value = "2016-03-09 04:36:26"
$dede.Range("D11").Value = "2016-03-09 04:36:26"
$dede.Range("D12").Value = "2016-03-09 04:36:26"
$dede.Range("D13").Value = "2016-03-09 04:36:26"
$dede.Range("D14").Value = "2016-03-09 04:36:26"
$dede.Range("D15").Value = "2016-03-09 04:36:26"
$dede.Range("D16").Value = "2016-03-09 04:36:26"
